$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Step 2) is rewritten / split: original single step becomes two
# rows of Steps/Expected Result content.
$ws.Range("C3").Value = "Step 2: Log in as a user with the appropriate role"
$ws.Range("D3").Value = "I am redirected to the user's dashboard"

$ws.Range("C4").Value = "Step 3: Go to the ""Team KPI"" page "
$ws.Range("D4").Value = "A list of KPIs is displayed "

$ws.Range("C5").Value = "Step 4: Delete one that belongs to someone on any of my teams"
$ws.Range("D5").Value = "The data is removed from the database."

$ws.Range("C6").Value = "Step 5: While logged in try to delete a kpi about me"
$ws.Range("D6").Value = "I am denied access to this"

$ws.Range("C7").Value = "Step 6: Try to delete a kpi of someone who is not under any of my teams"
$ws.Range("D7").Value = "I am denied access to this"

$ws.Range("C6:D7").WrapText = $true
$ws.Range("C6:D7").VerticalAlignment = -4160

$ws.Range("D3").Select()
